# Worked on Odd Even Sort: add a new day row (2/20/2013) to the workload log.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy the formatting of the last existing row (17) down into the new row (18)
# so the new row's cells pick up the same styles (date format, borders, alignment).
$ws.Range("A17:E17").Copy()
$ws.Range("A18").PasteSpecial(-4122)

# Fill in the new day's data.
$ws.Range("A18").Value = 41325
$ws.Range("B18").Value = "0H"
$ws.Range("C18").Value = "0.5H"
$ws.Range("D18").Value = "1H"
$ws.Range("E18").Value = "Odd-Even Sort"

# Match the author's final selection in the saved workbook.
$ws.Range("F18").Select()
